$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference style (no explicit fill/border/font overrides) used to avoid
# leaving a forced-text number format on cells after writing numeric-looking text
$plainStyle = $ws.Range("B2").Style

$ws.Range("D2").Value = "42.415.44"
$ws.Range("E2").Value = "  +0.32%  "
$ws.Range("D3").Value = "2.283.52"
$ws.Range("E3").Value = "  -0.09%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = $plainStyle
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("D5").Value = "'309.61"
$ws.Range("D5").Style = $plainStyle
$ws.Range("E5").Value = "  -4.17%  "
$ws.Range("D6").Value = "'103.02"
$ws.Range("D6").Style = $plainStyle
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("D7").Value = "'0.620"
$ws.Range("D7").Style = $plainStyle
$ws.Range("E7").Value = "  -1.44%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  -1.19%  "
$ws.Range("D10").Value = "'38.52"
$ws.Range("D10").Style = $plainStyle
$ws.Range("E10").Value = "  -2.98%  "
$ws.Range("D11").Value = "'0.0898"
$ws.Range("D11").Style = $plainStyle
$ws.Range("E11").Value = "  -0.84%  "
$ws.Range("D12").Value = "'8.18"
$ws.Range("D12").Style = $plainStyle
$ws.Range("E12").Value = "  -1.86%  "
$ws.Range("E13").Value = "  +0.53%  "
$ws.Range("D14").Value = "'0.968"
$ws.Range("D14").Style = $plainStyle
$ws.Range("E14").Value = "  -0.41%  "
$ws.Range("D15").Value = "'15.11"
$ws.Range("D15").Style = $plainStyle
$ws.Range("E15").Value = "  -0.06%  "
$ws.Range("D16").Value = "2.628.58"
$ws.Range("E16").Value = "  -0.16%  "
$ws.Range("D17").Value = "2.277.85"
$ws.Range("E17").Value = "  -0.09%  "
$ws.Range("D18").Value = "42.596.02"
$ws.Range("E18").Value = "  +0.54%  "
$ws.Range("E19").Value = "  -0.93%  "
$ws.Range("E20").Value = "  -1.30%  "
$ws.Range("E21").Value = "  +0.14%  "
$ws.Range("D22").Value = "'72.86"
$ws.Range("D22").Style = $plainStyle
$ws.Range("E22").Value = "  -0.31%  "
$ws.Range("D23").Value = "'3.38"
$ws.Range("D23").Style = $plainStyle
$ws.Range("E23").Value = "  -6.36%  "
$ws.Range("D24").Value = "'262.41"
$ws.Range("D24").Style = $plainStyle
$ws.Range("E24").Value = "  -2.02%  "
$ws.Range("E25").Value = "  -2.69%  "
$ws.Range("E26").Value = "  +0.32%  "
$ws.Range("D27").Value = "'10.66"
$ws.Range("D27").Style = $plainStyle
$ws.Range("E27").Value = "  -1.97%  "
$ws.Range("D28").Value = "'7.02"
$ws.Range("D28").Style = $plainStyle
$ws.Range("E28").Value = "  +14.39%  "
$ws.Range("D29").Value = "'2.33"
$ws.Range("D29").Style = $plainStyle
$ws.Range("E29").Value = "  +0.13%  "
$ws.Range("D30").Value = "'22.06"
$ws.Range("D30").Style = $plainStyle
$ws.Range("E30").Value = "  -1.91%  "
$ws.Range("D31").Value = "'35.71"
$ws.Range("D31").Style = $plainStyle
$ws.Range("E31").Value = "  -6.32%  "
$ws.Range("D32").Value = "'164.47"
$ws.Range("D32").Style = $plainStyle
$ws.Range("E32").Value = "  +0.12%  "
$ws.Range("E33").Value = "  -3.08%  "
$ws.Range("E34").Value = "  -3.19%  "
$ws.Range("E35").Value = "  +1.95%  "
$ws.Range("E36").Value = "  -3.11%  "
$ws.Range("E37").Value = "  -2.59%  "
$ws.Range("E38").Value = "  -2.61%  "
$ws.Range("D39").Value = "'2.73"
$ws.Range("D39").Style = $plainStyle
$ws.Range("E39").Value = "  -0.46%  "
$ws.Range("D40").Value = "'3.61"
$ws.Range("D40").Style = $plainStyle
$ws.Range("E40").Value = "  -2.73%  "
$ws.Range("E41").Value = "  +2.34%  "
$ws.Range("D42").Value = "'101.52"
$ws.Range("D42").Style = $plainStyle
$ws.Range("E42").Value = "  +11.29%  "
$ws.Range("D43").Value = "'69.08"
$ws.Range("D43").Style = $plainStyle
$ws.Range("E43").Value = "  -0.47%  "
$ws.Range("D44").Value = "'1.00"
$ws.Range("D44").Style = $plainStyle
$ws.Range("E44").Value = "  +0.22%  "
$ws.Range("D46").Value = "'11.98"
$ws.Range("D46").Style = $plainStyle
$ws.Range("E46").Value = "  -2.81%  "
$ws.Range("D47").Value = "1.716.00"
$ws.Range("E47").Value = "  +7.22%  "
$ws.Range("D48").Value = "'109.82"
$ws.Range("D48").Style = $plainStyle
$ws.Range("E48").Value = "  -2.56%  "
$ws.Range("D49").Value = "'77.14"
$ws.Range("D49").Style = $plainStyle
$ws.Range("E49").Value = "  -4.23%  "
$ws.Range("D50").Value = "'8.64"
$ws.Range("D50").Style = $plainStyle
$ws.Range("E50").Value = "  -2.88%  "
$ws.Range("D51").Value = "'5.13"
$ws.Range("D51").Style = $plainStyle
$ws.Range("E51").Value = "  -1.65%  "
